$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the capitalisation of three taxon labels in column A ("final main
# contribution" taxon-name cleanup): Dragonflies & Damselflies, Freshwater
# Crabs and Freshwater Fishes each had an inconsistently capitalised second
# word; update every occurrence in the data rows that use them (rows 9-14).
$ws.Range("A9").Value = "Dragonflies & damselflies"
$ws.Range("A10").Value = "Dragonflies & damselflies"
$ws.Range("A11").Value = "Freshwater crabs"
$ws.Range("A12").Value = "Freshwater crabs"
$ws.Range("A13").Value = "Freshwater fishes"
$ws.Range("A14").Value = "Freshwater fishes"

# Restore the sheet selection to a single cell (A10) as left by the author,
# instead of the prior multi-row selection/scroll position.
$ws.Activate()
$ws.Range("A10").Select()
